$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update site names (column A) to use spaces instead of underscores, proper case
$ws.Range("A2").Value = "Anacapa Passage"
$ws.Range("A3").Value = "Footprint"
$ws.Range("A4").Value = "Piggy Bank"

# Update location column (B) to site-specific slug values instead of the constant "anacapa_island"
$ws.Range("B2").Value = "anacapa_passage"
$ws.Range("B3").Value = "footprint"
$ws.Range("B4").Value = "piggy_bank"

# Update selection to A5 (matches the post-edit selection saved in the sheet)
$ws.Range("A5").Select()
